# Update the last row (row 41) of the schedule to reflect the extended
# "automated feature engineering" literature-review task:
#  - append a note about featuretools to the comments cell (I41)
#  - push the end-date/effort numbers out (F41, G41, H41)
#  - grow the row to fit the now-longer wrapped comment text
#  - leave the active selection on the edited cell, like the author did

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the featuretools sentence (on a new line) to the existing note in I41.
$existingNote = $ws.Range("I41").Value2
$addition = "הסבר טוב על featuretools https://towardsdatascience.com/automated-feature-engineering-in-python-99baf11cc219"
$ws.Range("I41").Value = $existingNote + "`n" + $addition

# Update the numeric work-tracking cells for the row.
$ws.Range("F41").Value = 44045.739583333336
$ws.Range("G41").Value = 1.5
$ws.Range("H41").Value = 0.75

# The longer wrapped comment needs a taller row.
$ws.Rows(41).RowHeight = 72.5

# Match the author's final selection/cursor position.
$ws.Range("I41").Select()
